# Apply the "Complete data collection for research" edit:
#  - Rename existing step sheets with numeric prefixes
#  - Insert a "-Results" leaderboard sheet after each step sheet
#  - Turn the 1-Prompt sheet's K1:O1 numeric labels into quoted-text labels
#  - Restore view state (selection/zoom/pane/active sheet) to match the
#    post-edit workbook
#
# NOTE: worksheet object references captured in a variable become stale
# (silently resolve to the wrong sheet) once any Worksheets.Add() call
# shifts sheet positions around them. So every sheet is re-looked-up by
# name ($wb.Worksheets.Item("...")) right before it is used, instead of
# being cached across Add() calls.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the four step sheets
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Prompts").Name = "1-Prompt"
$wb.Worksheets.Item("Step2").Name = "2-Adjective"
$wb.Worksheets.Item("Step3").Name = "3-Separator"
$wb.Worksheets.Item("Step4").Name = "4-Contamination"

# ---------------------------------------------------------------------------
# 2. Update the K1:O1 header cells on 1-Prompt: plain numbers -> quoted text
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("1-Prompt").Range("K1").Value = '"1"'
$wb.Worksheets.Item("1-Prompt").Range("L1").Value = '"2"'
$wb.Worksheets.Item("1-Prompt").Range("M1").Value = '"3"'
$wb.Worksheets.Item("1-Prompt").Range("N1").Value = '"4"'
$wb.Worksheets.Item("1-Prompt").Range("O1").Value = '"5"'

# ---------------------------------------------------------------------------
# 3. Insert the four "-Results" leaderboard sheets right after their
#    corresponding step sheet
# ---------------------------------------------------------------------------

function Style-Header($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Color = 16777215
    $rng.Interior.Color = 5855577
    $rng.HorizontalAlignment = -4108
}

function Fill-Results($sheetName, $rows) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Item(1, 1).Value = "Prompt"
    $ws.Cells.Item(1, 2).Value = "win"
    $ws.Cells.Item(1, 3).Value = "score"
    Style-Header $ws.Range("A1:C1")

    $r = 2
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $r = $r + 1
    }
}

# --- 1-Results (inserted after 1-Prompt) ---
$new1 = $wb.Worksheets.Add($null, $wb.Worksheets.Item("1-Prompt"))
$new1.Name = "1-Results"

# --- 2-Results (inserted after 2-Adjective) ---
$new2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item("2-Adjective"))
$new2.Name = "2-Results"

# --- 3-Results (inserted after 3-Separator) ---
$new3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item("3-Separator"))
$new3.Name = "3-Results"

# --- 4-Results (inserted after 4-Contamination) ---
$new4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item("4-Contamination"))
$new4.Name = "4-Results"

# Now fill in the data (each lookup is fresh, by name)
$rows1 = @(
    ,@('Reference:  ', 9, 362)
    ,@('Concept:    ', 11, 363)
    ,@('Info:       ', 7, 447)
    ,@('Information:', 8, 456)
    ,@('Context:    ', 5, 464)
    ,@('null        ', 8, 489)
    ,@('WebContext: ', 0, 530)
    ,@('Data:       ', 4, 541)
    ,@(2, 4, 546)
    ,@(5, 3, 569)
    ,@(4, 5, 609)
    ,@(3, 5, 629)
    ,@(1, 6, 656)
)
Fill-Results "1-Results" $rows1
$wb.Worksheets.Item("1-Results").Columns.Item(1).ColumnWidth = 9.140625

$rows2 = @(
    ,@('Useful concept:', 11, 413)
    ,@('Helpful concept:', 7, 532)
    ,@('Relevant concept:', 8, 532)
    ,@('Useful reference:', 5, 541)
    ,@('Referencial concept:', 5, 586)
    ,@('Reference:', 1, 620)
    ,@('Helpful reference:', 2, 630)
    ,@('Contextual concept:', 6, 635)
    ,@('Referencial reference:', 2, 635)
    ,@('Concept:', 3, 638)
    ,@('Verified concept:', 10, 651)
    ,@('Infomative concept: ', 4, 662)
    ,@('Infomative reference:', 2, 669)
    ,@('Relevant reference:', 4, 674)
    ,@('Contextual reference:', 2, 684)
    ,@('Verified reference:', 3, 743)
)
Fill-Results "2-Results" $rows2
$wb.Worksheets.Item("2-Results").Columns.Item(1).ColumnWidth = 20.28515625

$rows3 = @(
    ,@('std', 29, 196)
    ,@('(")', 16, 221)
    ,@('(*)', 11, 231)
    ,@("('')", 7, 258)
    ,@('#', 7, 306)
    ,@('###', 5, 337)
)
Fill-Results "3-Results" $rows3

$rows4 = @(
    ,@('std', 46, 118)
    ,@('<task>', 22, 139)
    ,@('<task> benchmark', 7, 186)
)
Fill-Results "4-Results" $rows4
$wb.Worksheets.Item("4-Results").Columns.Item(1).ColumnWidth = 20.5703125

# ---------------------------------------------------------------------------
# 4. Restore view state
# ---------------------------------------------------------------------------

# 1-Prompt: drop the forced zoom/tab-selected, select A1:C1 in the frozen pane
$wb.Worksheets.Item("1-Prompt").Activate()
$wb.Worksheets.Item("1-Prompt").Range("A1:C1").Select()
$excel.ActiveWindow.Zoom = 100

# 2-Adjective: scroll so column C becomes the frozen pane's top-left column,
# and land the selection on D39
$wb.Worksheets.Item("2-Adjective").Activate()
$wb.Worksheets.Item("2-Adjective").Range("D39").Select()
$excel.ActiveWindow.ScrollColumn = 3

# 1-Results: select the header row
$wb.Worksheets.Item("1-Results").Activate()
$wb.Worksheets.Item("1-Results").Range("A1:C1").Select()

# 2-Results: land just past the data
$wb.Worksheets.Item("2-Results").Activate()
$wb.Worksheets.Item("2-Results").Range("A19").Select()

# 4-Results: land on the last data row
$wb.Worksheets.Item("4-Results").Activate()
$wb.Worksheets.Item("4-Results").Range("A4").Select()

# 3-Results ends up the active sheet/tab
$wb.Worksheets.Item("3-Results").Activate()
$wb.Worksheets.Item("3-Results").Range("B5").Select()
